$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: quality_comparison
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# C1 / D1 are part of the merged header range B1:D1. We want them to end up
# with a plain (non-bold / non-aligned) style that only carries a border:
#   C1 -> border with top+bottom thin only      (matches styles.xml borderId=4)
#   D1 -> border with right+top+bottom thin     (matches styles.xml borderId=5)
# Resetting to "Normal" first drops the inherited bold/alignment formatting,
# then the borders are built starting from the full box border (already an
# existing border definition) and trimming edges off, which keeps both the
# border table and the cell style table free of unused/duplicate entries.
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1            # full box border
$c1.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft  -> none
$c1.Borders.Item(10).LineStyle = -4142  # xlEdgeRight -> none

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1            # full box border
$d1.Borders.Item(7).LineStyle = -4142   # xlEdgeLeft  -> none

# Rename header label from "fedcore" to "approach"
$ws1.Range("C2").Value = "approach"

# ---------------------------------------------------------------------------
# Sheet 2: computational_comparison
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# C1/D1 and F1/G1 need the exact same styles as sheet1's C1/D1 (two merged
# header ranges, B1:D1 and E1:G1). Copy the already-built formats instead of
# re-deriving them, so the shared style table stays minimal.
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
